$d = $word.ActiveDocument

$rsquo = [char]0x2019
$lsquo = [char]0x2018
$ldq = [char]0x201C
$rdq = [char]0x201D

function New-OOXML($bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Append a new run to the "Now I will give a brief introduction..."
#    paragraph (the "In Java" intro paragraph).
# ---------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(10)
$introPara.Range.InsertAfter(" Java has a properties class that can be used to implement this kind of idea. There" + $rsquo + "s a specific convention to use to create a properties file which can be read by included methods. This will then generate a kind of hash table that can be accessed easily.")

# ---------------------------------------------------------------------
# 2) Replace the "To access external files..." sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("To access external files in java you can use the BufferedReader class. Here is some example code:", $true, $false, $false, $false, $false, $true, 1, $false, "To access the external file to be read into the properties object:", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Replace the whole BufferedReader try/catch code sample (the 15
#    "Code" styled paragraphs right after the "Code:" label) with the
#    new Properties-based sample.
# ---------------------------------------------------------------------
$codeLabel = $d.Paragraphs.Item(13)
$firstCode = $d.Paragraphs.Item(14)
$lastCode = $d.Paragraphs.Item(28)
$codeRange = $d.Range($firstCode.Range.Start, $lastCode.Range.End)
$codeRange.Delete() | Out-Null

# After deletion, a new empty paragraph is inserted right after the
# "Code:" label (#13); its Range (the paragraph mark itself) is then
# expanded into the six replacement paragraphs via InsertXML, which
# lets us control run-splits, paragraph styles and empty paragraphs
# precisely.
$anchor = $d.Paragraphs.Item(13)
$anchor.Range.InsertParagraphAfter() | Out-Null
$newp = $d.Paragraphs.Item(14)

$newBlock = '<w:p><w:pPr><w:pStyle w:val="Code"/></w:pPr><w:r><w:t xml:space="preserve">        </w:t></w:r><w:r><w:t>Properties properties = new Properties();</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Code"/></w:pPr><w:r><w:t xml:space="preserve">        properties.load(' + $ldq + '&lt;root&gt;:\properties.txt' + $rdq + ');</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Code"/></w:pPr></w:p>' +
    '<w:p><w:r><w:t>You can then access the properties in this file using a Hash table-like key function:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Code"/></w:pPr><w:r><w:t>Code:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Code"/></w:pPr><w:r><w:t xml:space="preserve">       String thisPropertyIsSetTo = properties.get(thisProperty);</w:t></w:r></w:p>'

$newp.Range.InsertXML((New-OOXML $newBlock)) | Out-Null

# ---------------------------------------------------------------------
# 4) Replace the "This Code will allow you to read..." paragraph (the
#    bold-blue-heading-formatted paragraph right after the blank
#    paragraph that follows the code block), moving
#    <w:lastRenderedPageBreak/> onto its first run to match the new
#    pagination, and keeping its paragraph-mark run formatting intact.
# ---------------------------------------------------------------------
$full = $d.Content
$full.Find.Execute("This Code will allow you to read each line in a file and execute the required code on each line.") | Out-Null

$replacement = '<w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/></w:rPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">This Code will </w:t></w:r>' +
    '<w:r><w:t>find the property key ' + $lsquo + 'thisProperty' + $rsquo + ' and set the String thisPropertyIsSetTo to the value stored in the properties file as the value for ' + $lsquo + 'thisProperty' + $rsquo + '.</w:t></w:r>' +
    '</w:p>'

$full.InsertXML((New-OOXML $replacement)) | Out-Null
